$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprinklers")

# Sort the sprinkler data (A10:H58) by column A ascending, as recorded by
# Excel's Data > Sort dialog (records a <sortState> on the sheet too).
$rng = $ws.Range("A10:H58")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A10:A58"), 0, 1, $null, 0)
$ws.Sort.SetRange($rng)
$ws.Sort.Header = 0
$ws.Sort.Apply()

# Switch the active sheet/selection to "Sprinklers" (was "Test Parameters"),
# landing the cursor on A10.
$ws.Activate()
$ws.Range("A10").Select()
